# Update export mutasi saldo
# Insert a new "User" column between "Keterangan" (C) and "Kredit" (was D, now E)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D; this shifts the old D:F (Kredit/Debit/Saldo)
# to E:G and makes the new D inherit formatting from column C (text header style)
$ws.Columns("D:D").Insert() | Out-Null

# Populate the new header cell
$ws.Range("D1").Value() = "User"

# Match the new column's width to its numeric-column neighbours (E:G)
$ws.Columns("D:D").ColumnWidth() = 16.948043184885268

# Move/restore the active selection as in the authored workbook
$ws.Range("D2").Select() | Out-Null
